$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(263, 1029804860, "Update Laporan TS Harian", "2024-07-29 17:30:20.668248")
    ,@(264, 1029804860, "Kirim Laporan TS Harian", "2024-07-29 17:31:57.990398")
    ,@(265, 1029804860, "Update Laporan TS Harian", "2024-07-29 18:38:53.049810")
    ,@(266, 1029804860, "Kirim Laporan TS Harian", "2024-07-29 18:58:33.757631")
    ,@(267, 1029804860, "Add user", "2024-07-30 15:26:16.639377")
    ,@(268, -4146843555, "Update Laporan TS Harian", "2024-07-30 17:42:24.180796")
    ,@(269, -4146843555, "Update Laporan TS Harian", "2024-07-30 17:43:25.835575")
    ,@(270, -4146843555, "Kirim Laporan TS Harian", "2024-07-30 17:44:29.026010")
    ,@(271, -4146843555, "Update Laporan TS Harian", "2024-07-30 17:45:08.311373")
    ,@(272, -4146843555, "Update Laporan TS Harian", "2024-07-30 17:45:33.992416")
    ,@(273, -4146843555, "Update Laporan TS Harian", "2024-07-30 17:46:16.152798")
    ,@(274, -4146843555, "Update Laporan TS Harian", "2024-07-30 17:46:36.125724")
    ,@(275, 1029804860, "Update Laporan TS Harian", "2024-07-30 18:20:13.367459")
    ,@(276, 1029804860, "Kirim Laporan TS Harian", "2024-07-30 18:20:46.156950")
    ,@(277, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-30 18:34:39.703646")
    ,@(278, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-30 18:36:50.409450")
    ,@(279, -4146843555, "Kirim Laporan TS Harian", "2024-07-30 18:37:46.339847")
    ,@(280, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-30 18:51:37.265321")
    ,@(281, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-30 18:58:09.744323")
    ,@(282, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Otomatis", "2024-07-30 19:02:08.416197")
    ,@(283, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-30 19:15:50.316209")
    ,@(284, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-30 19:27:05.002617")
    ,@(285, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-30 19:28:41.808572")
    ,@(286, 1029804860, "Update Laporan TS Harian", "2024-07-31 18:21:24.934501")
    ,@(287, 1029804860, "Update Laporan TS Harian", "2024-07-31 18:37:28.130751")
    ,@(288, 1029804860, "Kirim Laporan TS Harian", "2024-07-31 18:45:17.702898")
    ,@(289, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Manual", "2024-07-31 18:57:12.414347")
    ,@(290, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Otomatis", "2024-07-31 19:16:16.943330")
    ,@(291, -4146843555, "Mengirim Laporan Realisasi Tindak Lanjut JN Max Otomatis", "2024-07-31 19:20:38.577126")
    ,@(292, -4146843555, "Kirim Laporan TS Harian", "2024-07-31 20:22:46.185023")
)

$startRow = 264
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $item = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $item[0]
    $ws.Cells.Item($r, 2).Value = $item[1]
    $ws.Cells.Item($r, 3).Value = $item[2]
    $ws.Cells.Item($r, 4).Value = $item[3]
}

Write-Output "Added $($rows.Count) rows"